# Apply the "hello" -> "Hello" rename and add DISASTER events / new EU values
# for a handful of log rows in the "Output Schedules" sheet, per the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Output Schedules")

# Row 5 = Turn #1: hello -> Turn #1: Hello, with an updated TRANSFORM action/EU.
$ws.Range("A5").Value = "Turn #1: Hello: "
$ws.Range("B5").Value = "(TRANSFORM Hello (INPUTS (Population 5) (MetallicElements 3) (MetallicAlloys 2)) (OUTPUTS (Population 5) (Electronics 2) (ElectonicsWaste 1))) EU: 0.22096544749376973"

# Row 7 = Turn #2: Erewhon gets hit by a DISASTER (Earthquake) before its TRANSFORM.
$ws.Range("B7").Value = "(DISASTER Erewhon (Earthquake)) (TRANSFORM Erewhon (INPUTS (Population 5) (MetallicElements 3) (MetallicAlloys 2)) (OUTPUTS (Population 5) (Electronics 2) (ElectonicsWaste 1))) EU: -2.0130409724967318"

# Row 9 = Turn #2: hello -> Turn #2: Hello, hit by a DISASTER (Fire) before its TRANSFORM.
$ws.Range("A9").Value = "Turn #2: Hello: "
$ws.Range("B9").Value = "(DISASTER Hello (Fire)) (TRANSFORM Hello (INPUTS (Population 1) (MetallicElements 2))  (OUTPUTS (Population 1) (MetallicAlloys 1) (MetallicAlloysWaste 1))) EU: 0.16863274528787145"

# Row 11 = Turn #3: Erewhon gets hit by a DISASTER (Tornado) and simply PASSES.
$ws.Range("B11").Value = "(DISASTER Erewhon (Tornado)) (PASSES Erewhon )"

# Row 13 = Turn #3: hello -> Turn #3: Hello, with an updated TRANSFORM action/EU.
$ws.Range("A13").Value = "Turn #3: Hello: "
$ws.Range("B13").Value = "(TRANSFORM Hello (INPUTS (Population 5) (MetallicElements 3) (MetallicAlloys 2)) (OUTPUTS (Population 5) (Electronics 2) (ElectonicsWaste 1))) EU: 0.19193310620438808"
